# Rewrite the three worksheets ("Blok Sensus", "SLS", "DESA") with the new
# sample header/data rows, matching the updated template_import_wilkerstat.xlsx
# contents described in the commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Blok Sensus" -> Kode/Nama Blok Sensus, numeric codes
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Blok Sensus")

$ws1.Range("A1").Value = "Kode Blok Sensus"
$ws1.Range("B1").Value = "Nama Blok Sensus"
$ws1.Range("A1:B1").Font.Bold = $true

$ws1.Range("A2").Value = 7371100001001
$ws1.Range("B2").Value = "BS 001 DESA CONTOH"

$ws1.Range("A3").Value = 7371100001002
$ws1.Range("B3").Value = "BS 002 DESA CONTOH"

# ---------------------------------------------------------------------------
# Sheet "SLS" -> Kode/Nama SLS, text codes
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("SLS")

$ws2.Range("A1").Value = "Kode SLS"
$ws2.Range("B1").Value = "Nama SLS"
$ws2.Range("A1:B1").Font.Bold = $true

# Keep the leading-zero codes as literal text (matches the source export's
# t="s" shared-string cells) instead of letting them auto-coerce to numbers.
$ws2.Range("A2").NumberFormat = "@"
$ws2.Range("A2").Value = "001"
$ws2.Range("B2").Value = "SLS 001 DESA CONTOH"

$ws2.Range("A3").NumberFormat = "@"
$ws2.Range("A3").Value = "002"
$ws2.Range("B3").Value = "SLS 002 DESA CONTOH"

# ---------------------------------------------------------------------------
# Sheet "DESA" -> Kode/Nama Desa, numeric codes
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("DESA")

$ws3.Range("A1").Value = "Kode Desa"
$ws3.Range("B1").Value = "Nama Desa"
$ws3.Range("A1:B1").Font.Bold = $true

$ws3.Range("A2").Value = 7371100001
$ws3.Range("B2").Value = "DESA CONTOH 1"

$ws3.Range("A3").Value = 7371100002
$ws3.Range("B3").Value = "DESA CONTOH 2"

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping: the updated workbook activates the
# first sheet ("Blok Sensus") instead of "DESA", and every sheet's selection
# resets to its header row (A1:B1).
# ---------------------------------------------------------------------------
$ws2.Range("A1:B1").Select() | Out-Null
$ws3.Range("A1:B1").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A1:B1").Select() | Out-Null
